$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Range("B1").Value = "NUMERO"
$ws.Range("C1").Value = "valor"
$ws.Range("D1").Value = "vencimento"
$ws.Range("E1").Value = "status"
$ws.Range("F1").Value = "STATUS"

# --- Clear old data rows (2-6) then rewrite rows 2-5 ---
$ws.Range("A2:F6").ClearContents()

# Row 2 - Gustavo
$ws.Range("A2").Value = "Gustavo"
$ws.Range("B2").Value = 619994121708
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 45910
$ws.Range("D2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = "pendente"
$ws.Range("F2").Value = "enviado"

# Row 3 - Gisele
$ws.Range("A3").Value = "Gisele"
$ws.Range("B3").Value = 61993367127
$ws.Range("C3").Value = 85
$ws.Range("D3").Value = 45912
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E3").Value = "pendente"
$ws.Range("F3").Value = "enviado"

# Row 4 - Iury
$ws.Range("A4").Value = "Iury"
$ws.Range("B4").Value = 61993121599
$ws.Range("C4").Value = 90
$ws.Range("D4").Value = 45915
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E4").Value = "pendente"
$ws.Range("F4").Value = "enviado"

# Row 5 - Murillo
$ws.Range("A5").Value = "Murillo"
$ws.Range("B5").Value = 61996544168
$ws.Range("C5").Value = 75
$ws.Range("D5").Value = 45912
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = "pendente"
$ws.Range("F5").Value = "enviado"

# Remove the now-unused row 6 (shrinks used range to A1:F5)
$ws.Rows("6:6").Delete()
